$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new product data
$ws.Range("A2").Value = "Esencia"
$ws.Range("B2").Value = 400
$ws.Range("C2").Value = "gramos"

$ws.Range("A3").Value = "Bolsa de Regalo"
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = "unidades"

$ws.Range("A4").Value = "Splash Pink"
$ws.Range("B4").Value = 28
$ws.Range("C4").Value = "unidades"

# Delete rows 5 through 11 (now obsolete products)
$ws.Range("A5:C11").EntireRow.Delete()
